$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1023.207113
$ws.Range("H2").Value = 3069.621339
$ws.Range("I2").Value = 0.206942302533768
$ws.Range("J2").Value = 0.223591263450141
$ws.Range("M2").Value = 16.14072933333334
$ws.Range("N2").Value = 48.42218800000001
$ws.Range("O2").Value = 0.03423048004954622
$ws.Range("P2").Value = 0.03634868370049611
$ws.Range("Q2").Value = 16515.30906287442
$ws.Range("R2").Value = 148637.7815658698
$ws.Range("S2").Value = 0.007083734358289305
$ws.Range("T2").Value = 0.008127248113343472
$ws.Range("G3").Value = 1023.207113
$ws.Range("H3").Value = 3069.621339
$ws.Range("I3").Value = 0.206942302533768
$ws.Range("J3").Value = 0.223591263450141
$ws.Range("O3").Value = 0.1719151703242873
$ws.Range("P3").Value = 0.1825533892714798
$ws.Range("Q3").Value = 82944.56187563574
$ws.Range("R3").Value = 746501.0568807218
$ws.Range("S3").Value = 0.03557652118739293
$ws.Range("T3").Value = 0.04081734295431559
$ws.Range("G4").Value = 1023.207113
$ws.Range("H4").Value = 3069.621339
$ws.Range("I4").Value = 0.206942302533768
$ws.Range("J4").Value = 0.223591263450141
$ws.Range("M4").Value = 168.70371
$ws.Range("N4").Value = 506.11113
$ws.Range("O4").Value = 0.3577786889414888
$ws.Range("P4").Value = 0.3799182594076638
$ws.Range("Q4").Value = 172618.8360614892
$ws.Range("R4").Value = 1553569.524553403
$ws.Range("S4").Value = 0.07403954568706446
$ws.Range("T4").Value = 0.08494640362873797
$ws.Range("G5").Value = 1023.207113
$ws.Range("H5").Value = 3069.621339
$ws.Range("I5").Value = 0.206942302533768
$ws.Range("J5").Value = 0.223591263450141
$ws.Range("M5").Value = 82.43477250000001
$ws.Range("N5").Value = 164.869545
$ws.Range("O5").Value = 0.1748236883957081
$ws.Range("P5").Value = 0.1237612588479007
$ws.Range("Q5").Value = 84347.8455805368
$ws.Range("R5").Value = 506087.0734832208
$ws.Range("S5").Value = 0.03617841661405381
$ws.Range("T5").Value = 0.02767193623198206
$ws.Range("G6").Value = 1023.207113
$ws.Range("H6").Value = 3069.621339
$ws.Range("I6").Value = 0.206942302533768
$ws.Range("J6").Value = 0.223591263450141
$ws.Range("M6").Value = 123.1883796666667
$ws.Range("N6").Value = 369.565139
$ws.Range("O6").Value = 0.2612519722889696
$ws.Range("P6").Value = 0.2774184087724594
$ws.Range("Q6").Value = 126047.2263138779
$ws.Range("R6").Value = 1134425.036824901
$ws.Range("S6").Value = 0.05406408468696752
$ws.Range("T6").Value = 0.06202833252176188
$ws.Range("I7").Value = 0.2556293015703424
$ws.Range("J7").Value = 0.2761952380599582
$ws.Range("M7").Value = 16.14072933333334
$ws.Range("N7").Value = 48.42218800000001
$ws.Range("O7").Value = 0.03423048004954622
$ws.Range("P7").Value = 0.03634868370049611
$ws.Range("Q7").Value = 20400.84056894088
$ws.Range("R7").Value = 183607.5651204679
$ws.Range("S7").Value = 0.008750313707483038
$ws.Range("T7").Value = 0.01003933334782465
$ws.Range("I8").Value = 0.2556293015703424
$ws.Range("J8").Value = 0.2761952380599582
$ws.Range("O8").Value = 0.1719151703242873
$ws.Range("P8").Value = 0.1825533892714798
$ws.Range("S8").Value = 0.04394655491934402
$ws.Range("T8").Value = 0.05042037680848859
$ws.Range("I9").Value = 0.2556293015703424
$ws.Range("J9").Value = 0.2761952380599582
$ws.Range("M9").Value = 168.70371
$ws.Range("N9").Value = 506.11113
$ws.Range("O9").Value = 0.3577786889414888
$ws.Range("P9").Value = 0.3799182594076638
$ws.Range("Q9").Value = 213230.6056326185
$ws.Range("R9").Value = 1919075.450693566
$ws.Range("S9").Value = 0.09145871637086556
$ws.Range("T9").Value = 0.1049316141004246
$ws.Range("I10").Value = 0.2556293015703424
$ws.Range("J10").Value = 0.2761952380599582
$ws.Range("M10").Value = 82.43477250000001
$ws.Range("N10").Value = 164.869545
$ws.Range("O10").Value = 0.1748236883957081
$ws.Range("P10").Value = 0.1237612588479007
$ws.Range("Q10").Value = 104192.2342156087
$ws.Range("R10").Value = 625153.4052936521
$ws.Range("S10").Value = 0.04469005736254602
$ws.Range("T10").Value = 0.03418227035009604
$ws.Range("I11").Value = 0.2556293015703424
$ws.Range("J11").Value = 0.2761952380599582
$ws.Range("M11").Value = 123.1883796666667
$ws.Range("N11").Value = 369.565139
$ws.Range("O11").Value = 0.2612519722889696
$ws.Range("P11").Value = 0.2774184087724594
$ws.Range("Q11").Value = 155702.1644824781
$ws.Range("R11").Value = 1401319.480342302
$ws.Range("S11").Value = 0.06678365921010372
$ws.Range("T11").Value = 0.07662164345312422
$ws.Range("G12").Value = 761.6879476666667
$ws.Range("H12").Value = 2285.063843
$ws.Range("I12").Value = 0.1540503928282995
$ws.Range("J12").Value = 0.1664440839100528
$ws.Range("M12").Value = 16.14072933333334
$ws.Range("N12").Value = 48.42218800000001
$ws.Range("O12").Value = 0.03423048004954622
$ws.Range("P12").Value = 0.03634868370049611
$ws.Range("Q12").Value = 12294.19899974983
$ws.Range("R12").Value = 110647.7909977485
$ws.Range("S12").Value = 0.005273218898333864
$ws.Range("T12").Value = 0.006050023359865342
$ws.Range("G13").Value = 761.6879476666667
$ws.Range("H13").Value = 2285.063843
$ws.Range("I13").Value = 0.1540503928282995
$ws.Range("J13").Value = 0.1664440839100528
$ws.Range("O13").Value = 0.1719151703242873
$ws.Range("P13").Value = 0.1825533892714798
$ws.Range("Q13").Value = 61744.95104899044
$ws.Range("R13").Value = 555704.559440914
$ws.Range("S13").Value = 0.02648359952160047
$ws.Range("T13").Value = 0.03038493164196672
$ws.Range("G14").Value = 761.6879476666667
$ws.Range("H14").Value = 2285.063843
$ws.Range("I14").Value = 0.1540503928282995
$ws.Range("J14").Value = 0.1664440839100528
$ws.Range("M14").Value = 168.70371
$ws.Range("N14").Value = 506.11113
$ws.Range("O14").Value = 0.3577786889414888
$ws.Range("P14").Value = 0.3799182594076638
$ws.Range("Q14").Value = 128499.5826336525
$ws.Range("R14").Value = 1156496.243702873
$ws.Range("S14").Value = 0.05511594757703031
$ws.Range("T14").Value = 0.06323514664781039
$ws.Range("G15").Value = 761.6879476666667
$ws.Range("H15").Value = 2285.063843
$ws.Range("I15").Value = 0.1540503928282995
$ws.Range("J15").Value = 0.1664440839100528
$ws.Range("M15").Value = 82.43477250000001
$ws.Range("N15").Value = 164.869545
$ws.Range("O15").Value = 0.1748236883957081
$ws.Range("P15").Value = 0.1237612588479007
$ws.Range("Q15").Value = 62789.57268189358
$ws.Range("R15").Value = 376737.4360913615
$ws.Range("S15").Value = 0.02693165787305105
$ws.Range("T15").Value = 0.02059932935249375
$ws.Range("G16").Value = 761.6879476666667
$ws.Range("H16").Value = 2285.063843
$ws.Range("I16").Value = 0.1540503928282995
$ws.Range("J16").Value = 0.1664440839100528
$ws.Range("M16").Value = 123.1883796666667
$ws.Range("N16").Value = 369.565139
$ws.Range("O16").Value = 0.2612519722889696
$ws.Range("P16").Value = 0.2774184087724594
$ws.Range("Q16").Value = 93831.10408468546
$ws.Range("R16").Value = 844479.9367621691
$ws.Range("S16").Value = 0.04024596895828377
$ws.Range("T16").Value = 0.04617465290791656
$ws.Range("G17").Value = 1104.505371
$ws.Range("H17").Value = 2209.010742
$ws.Range("I17").Value = 0.2233847690576539
$ws.Range("J17").Value = 0.1609043749153822
$ws.Range("M17").Value = 16.14072933333334
$ws.Range("N17").Value = 48.42218800000001
$ws.Range("O17").Value = 0.03423048004954622
$ws.Range("P17").Value = 0.03634868370049611
$ws.Range("Q17").Value = 17827.52224052392
$ws.Range("R17").Value = 106965.1334431435
$ws.Range("S17").Value = 0.007646567880600511
$ws.Range("T17").Value = 0.005848662229825268
$ws.Range("G18").Value = 1104.505371
$ws.Range("H18").Value = 2209.010742
$ws.Range("I18").Value = 0.2233847690576539
$ws.Range("J18").Value = 0.1609043749153822
$ws.Range("O18").Value = 0.1719151703242873
$ws.Range("P18").Value = 0.1825533892714798
$ws.Range("Q18").Value = 89534.86828123871
$ws.Range("R18").Value = 537209.2096874323
$ws.Range("S18").Value = 0.03840323062039815
$ws.Range("T18").Value = 0.0293736389894119
$ws.Range("G19").Value = 1104.505371
$ws.Range("H19").Value = 2209.010742
$ws.Range("I19").Value = 0.2233847690576539
$ws.Range("J19").Value = 0.1609043749153822
$ws.Range("M19").Value = 168.70371
$ws.Range("N19").Value = 506.11113
$ws.Range("O19").Value = 0.3577786889414888
$ws.Range("P19").Value = 0.3799182594076638
$ws.Range("Q19").Value = 186334.1538026264
$ws.Range("R19").Value = 1118004.922815759
$ws.Range("S19").Value = 0.07992230980294465
$ws.Range("T19").Value = 0.06113051004893017
$ws.Range("G20").Value = 1104.505371
$ws.Range("H20").Value = 2209.010742
$ws.Range("I20").Value = 0.2233847690576539
$ws.Range("J20").Value = 0.1609043749153822
$ws.Range("M20").Value = 82.43477250000001
$ws.Range("N20").Value = 164.869545
$ws.Range("O20").Value = 0.1748236883957081
$ws.Range("P20").Value = 0.1237612588479007
$ws.Range("Q20").Value = 91049.64898341312
$ws.Range("R20").Value = 364198.5959336525
$ws.Range("S20").Value = 0.03905294925808248
$ws.Range("T20").Value = 0.01991372799366228
$ws.Range("G21").Value = 1104.505371
$ws.Range("H21").Value = 2209.010742
$ws.Range("I21").Value = 0.2233847690576539
$ws.Range("J21").Value = 0.1609043749153822
$ws.Range("M21").Value = 123.1883796666667
$ws.Range("N21").Value = 369.565139
$ws.Range("O21").Value = 0.2612519722889696
$ws.Range("P21").Value = 0.2774184087724594
$ws.Range("Q21").Value = 136062.2269866205
$ws.Range("R21").Value = 816373.3619197232
$ws.Range("S21").Value = 0.05835971149562805
$ws.Range("T21").Value = 0.04463783565355257
$ws.Range("G22").Value = 791.0717773333332
$ws.Range("H22").Value = 2373.215332
$ws.Range("I22").Value = 0.1599932340099362
$ws.Range("J22").Value = 0.1728650396644658
$ws.Range("M22").Value = 16.14072933333334
$ws.Range("N22").Value = 48.42218800000001
$ws.Range("O22").Value = 0.03423048004954622
$ws.Range("P22").Value = 0.03634868370049611
$ws.Range("Q22").Value = 12768.47544117627
$ws.Range("R22").Value = 114916.2789705864
$ws.Range("S22").Value = 0.0054766452048395
$ws.Range("T22").Value = 0.006283416649637384
$ws.Range("G23").Value = 791.0717773333332
$ws.Range("H23").Value = 2373.215332
$ws.Range("I23").Value = 0.1599932340099362
$ws.Range("J23").Value = 0.1728650396644658
$ws.Range("O23").Value = 0.1719151703242873
$ws.Range("P23").Value = 0.1825533892714798
$ws.Range("Q23").Value = 64126.90172834422
$ws.Range("R23").Value = 577142.115555098
$ws.Range("S23").Value = 0.02750526407555174
$ws.Range("T23").Value = 0.03155709887729704
$ws.Range("G24").Value = 791.0717773333332
$ws.Range("H24").Value = 2373.215332
$ws.Range("I24").Value = 0.1599932340099362
$ws.Range("J24").Value = 0.1728650396644658
$ws.Range("M24").Value = 168.70371
$ws.Range("N24").Value = 506.11113
$ws.Range("O24").Value = 0.3577786889414888
$ws.Range("P24").Value = 0.3799182594076638
$ws.Range("Q24").Value = 133456.7437124272
$ws.Range("R24").Value = 1201110.693411845
$ws.Range("S24").Value = 0.05724216950358378
$ws.Range("T24").Value = 0.06567458498176063
$ws.Range("G25").Value = 791.0717773333332
$ws.Range("H25").Value = 2373.215332
$ws.Range("I25").Value = 0.1599932340099362
$ws.Range("J25").Value = 0.1728650396644658
$ws.Range("M25").Value = 82.43477250000001
$ws.Range("N25").Value = 164.869545
$ws.Range("O25").Value = 0.1748236883957081
$ws.Range("P25").Value = 0.1237612588479007
$ws.Range("Q25").Value = 65211.82199564399
$ws.Range("R25").Value = 391270.9319738639
$ws.Range("S25").Value = 0.02797060728797468
$ws.Range("T25").Value = 0.02139399491966658
$ws.Range("G26").Value = 791.0717773333332
$ws.Range("H26").Value = 2373.215332
$ws.Range("I26").Value = 0.1599932340099362
$ws.Range("J26").Value = 0.1728650396644658
$ws.Range("M26").Value = 123.1883796666667
$ws.Range("N26").Value = 369.565139
$ws.Range("O26").Value = 0.2612519722889696
$ws.Range("P26").Value = 0.2774184087724594
$ws.Range("Q26").Value = 97450.85044972344
$ws.Range("R26").Value = 877057.6540475111
$ws.Range("S26").Value = 0.04179854793798647
$ws.Range("T26").Value = 0.0479559442361042
